$wb = $excel.ActiveWorkbook

# Sheet 1: "展览" (first worksheet) - values updated
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F6").Value = 1099
$ws1.Range("F8").Value = 1794
$ws1.Range("F9").Value = 6427
$ws1.Range("F13").Value = 105
$ws1.Range("F14").Value = 379
$ws1.Range("F16").Value = 6749
$ws1.Range("F17").Value = 278
$ws1.Range("F18").Value = 1297
$ws1.Range("F21").Value = 222
$ws1.Range("F25").Value = 160
$ws1.Range("F27").Value = 101
$ws1.Range("F29").Value = 395
$ws1.Range("F30").Value = 434

# Sheet 4: "全部类型" (fourth worksheet) - values updated
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F6").Value = 1099
$ws4.Range("F8").Value = 1794
$ws4.Range("F9").Value = 6427
$ws4.Range("F13").Value = 105
$ws4.Range("F14").Value = 379
$ws4.Range("F16").Value = 6750
$ws4.Range("F17").Value = 278
$ws4.Range("F18").Value = 1297
$ws4.Range("F21").Value = 222
$ws4.Range("F25").Value = 160
$ws4.Range("F27").Value = 101
$ws4.Range("F29").Value = 395
$ws4.Range("F30").Value = 435
